{"js": "const table = context.document.body.tables.getFirst();\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst newValues = [\"87-54=\", \"35-26=\", \"31+46=\", \"2+85=\", \"17+31=\", \"44-0=\", \"17+46=\", \"21+16=\", \"92-36=\", \"78+9=\", \"62+12=\", \"49+39=\", \"87-16=\", \"31-30=\", \"10+40=\", \"97-49=\", \"8+20=\", \"90-63=\", \"25+71=\", \"61-34=\", \"99-45=\", \"73-18=\", \"85+7=\", \"0+25=\", \"67-66=\", \"92-12=\", \"88-65=\", \"46-22=\", \"50-41=\", \"91-42=\", \"21+45=\", \"82-47=\", \"6+81=\", \"93-91=\", \"63-42=\", \"91-30=\", \"14+27=\", \"30+0=\", \"44-27=\", \"15+73=\", \"77-71=\", \"30+63=\", \"90-24=\", \"37-26=\", \"88+11=\", \"50-5=\", \"50-34=\", \"60+4=\", \"59+10=\", \"80-35=\", \"80-15=\", \"28+35=\", \"58+33=\", \"14+21=\", \"66-41=\", \"58-16=\", \"63-27=\", \"72-11=\", \"43-7=\", \"97-32=\", \"59-30=\", \"14-14=\", \"76-40=\", \"7+73=\", \"67-29=\", \"25-12=\", \"59-7=\", \"83-15=\", \"93-55=\", \"32+48=\", \"37+2=\", \"30-2=\", \"7+2=\", \"53-47=\", \"33-14=\", \"50+48=\", \"78+4=\", \"51+28=\", \"65+24=\", \"71-41=\", \"69-61=\", \"24-3=\", \"17+28=\", \"53-45=\", \"54-25=\", \"35+0=\", \"35+58=\", \"42-38=\", \"26+68=\", \"4+57=\", \"26+59=\", \"94-7=\", \"27-9=\", \"50+15=\", \"58-51=\", \"96-3=\", \"9+62=\", \"76+23=\", \"99-16=\", \"18+8=\"];\n\nconst rowCount = table.rowCount;\nconst colCount = Math.round(newValues.length / rowCount);\n\nlet idx = 0;\nfor (let r = 0; r < rowCount; r++) {\n  for (let c = 0; c < colCount; c++) {\n    table.getCell(r, c).value = newValues[idx];\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$newValues = @(\n    \"87-54=\",\n    \"35-26=\",\n    \"31+46=\",\n    \"2+85=\",\n    \"17+31=\",\n    \"44-0=\",\n    \"17+46=\",\n    \"21+16=\",\n    \"92-36=\",\n    \"78+9=\",\n    \"62+12=\",\n    \"49+39=\",\n    \"87-16=\",\n    \"31-30=\",\n    \"10+40=\",\n    \"97-49=\",\n    \"8+20=\",\n    \"90-63=\",\n    \"25+71=\",\n    \"61-34=\",\n    \"99-45=\",\n    \"73-18=\",\n    \"85+7=\",\n    \"0+25=\",\n    \"67-66=\",\n    \"92-12=\",\n    \"88-65=\",\n    \"46-22=\",\n    \"50-41=\",\n    \"91-42=\",\n    \"21+45=\",\n    \"82-47=\",\n    \"6+81=\",\n    \"93-91=\",\n    \"63-42=\",\n    \"91-30=\",\n    \"14+27=\",\n    \"30+0=\",\n    \"44-27=\",\n    \"15+73=\",\n    \"77-71=\",\n    \"30+63=\",\n    \"90-24=\",\n    \"37-26=\",\n    \"88+11=\",\n    \"50-5=\",\n    \"50-34=\",\n    \"60+4=\",\n    \"59+10=\",\n    \"80-35=\",\n    \"80-15=\",\n    \"28+35=\",\n    \"58+33=\",\n    \"14+21=\",\n    \"66-41=\",\n    \"58-16=\",\n    \"63-27=\",\n    \"72-11=\",\n    \"43-7=\",\n    \"97-32=\",\n    \"59-30=\",\n    \"14-14=\",\n    \"76-40=\",\n    \"7+73=\",\n    \"67-29=\",\n    \"25-12=\",\n    \"59-7=\",\n    \"83-15=\",\n    \"93-55=\",\n    \"32+48=\",\n    \"37+2=\",\n    \"30-2=\",\n    \"7+2=\",\n    \"53-47=\",\n    \"33-14=\",\n    \"50+48=\",\n    \"78+4=\",\n    \"51+28=\",\n    \"65+24=\",\n    \"71-41=\",\n    \"69-61=\",\n    \"24-3=\",\n    \"17+28=\",\n    \"53-45=\",\n    \"54-25=\",\n    \"35+0=\",\n    \"35+58=\",\n    \"42-38=\",\n    \"26+68=\",\n    \"4+57=\",\n    \"26+59=\",\n    \"94-7=\",\n    \"27-9=\",\n    \"50+15=\",\n    \"58-51=\",\n    \"96-3=\",\n    \"9+62=\",\n    \"76+23=\",\n    \"99-16=\",\n    \"18+8=\"\n)\n\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n    for ($c = 1; $c -le $cols; $c++) {\n        $t.Cell($r, $c).Range.Text = $newValues[$idx]\n        $idx++\n    }\n}\n"}
